$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add P1 = 14, Q1 = 15, with the same style as the other header cells (s="1")
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# For data rows 2-25: add new columns P and Q (both value 2),
# and flip the values in columns I, K, M, O (swap 1<->2)
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1 (was 2)
    $ws.Cells.Item($r, 16).Value = 2  # P = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2
}
